$d = $word.ActiveDocument

# --- 1. Insert a new task bullet at the top of the "Tasks undertaken" list ---
# Locate the paragraph that currently reads "Develop message board".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Develop message board") {
        $targetIndex = $i
        break
    }
}

# Insert a brand-new (empty) list paragraph immediately before it - it
# inherits the same pPr / list numbering / rStyle automatically.
$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphBefore()

# Re-fetch by index (the old $target reference now points at the new, empty
# paragraph rather than "Develop message board", which shifted down by one).
$newPara = $d.Paragraphs.Item($targetIndex)

$newText = "Use ComponentTest.json to test out the messaging system as it is the one with the components done"
$newPara.Range.Text = $newText

# Colour just the run's text (not the paragraph mark) red, matching the diff.
$start = $newPara.Range.Start
$textRange = $d.Range($start, $start + $newText.Length)
$textRange.Font.Color = 255

# --- 2. Clear the text of the last "What we found out" bullet, leaving the
#        empty paragraph (with its pPr) in place. ---
$foundIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Describe the outcomes, and how they relate to the spike topic + graphs/screenshots/outputs as needed") {
        $foundIndex = $i
        break
    }
}

$found = $d.Paragraphs.Item($foundIndex)
$fullText = $found.Range.Text
$trimmedLen = $fullText.Length - 1
$clearRange = $d.Range($found.Range.Start, $found.Range.Start + $trimmedLen)
$clearRange.Delete()
